# Split the single "UK_utility_adjustment" time-series-factor sheet into
# three activity-status-specific sheets (single males, single females,
# couples), per the "Working activity status alignment" commit.
#
# 1) Rename the existing sheet to the "single males" variant.
# 2) Duplicate it twice (same data, still to be updated from the real
#    targets later per the commit message) for "single females" and
#    "couples", inserting each copy right after its source so the three
#    end up adjacent, in the order smales -> sfemales -> couples.
# 3) Leave the last of the three ("couples") as the active/selected sheet,
#    matching the workbook's saved tab selection.

$wb = $excel.ActiveWorkbook

$smales = $wb.Worksheets.Item("UK_utility_adjustment")
$smales.Name = "UK_utility_adj_smales"

$smales.Copy([System.Reflection.Missing]::Value, $smales)
$sfemales = $wb.Worksheets.Item("UK_utility_adj_smales (2)")
$sfemales.Name = "UK_utility_adj_sfemales"

$sfemales.Copy([System.Reflection.Missing]::Value, $sfemales)
$couples = $wb.Worksheets.Item("UK_utility_adj_sfemales (2)")
$couples.Name = "UK_utility_adj_couples"

$couples.Activate()
